$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-classify destination purposes (summarization script update):
# collapse "Health and Exercise" + "Social/Recreation/Eat Meal" into "Recreation/Eat Meal"
$ws.Range("B11").Value = "Recreation/Eat Meal"
$ws.Range("B13").Value = "Recreation/Eat Meal"
$ws.Range("B14").Value = "Recreation/Eat Meal"
$ws.Range("B15").Value = "Recreation/Eat Meal"
$ws.Range("B16").Value = "Recreation/Eat Meal"
$ws.Range("B17").Value = "Recreation/Eat Meal"
$ws.Range("B20").Value = "Recreation/Eat Meal"

# collapse "Shop" + "Errands" (+ others) into "Errands and Shopping"
$ws.Range("B4").Value = "Errands and Shopping"
$ws.Range("B8").Value = "Errands and Shopping"
$ws.Range("B9").Value = "Errands and Shopping"
$ws.Range("B10").Value = "Errands and Shopping"
$ws.Range("B12").Value = "Errands and Shopping"
$ws.Range("B19").Value = "Errands and Shopping"
$ws.Range("B21").Value = "Errands and Shopping"

# "Other" category (transfer to another mode) no longer has a mapped value
$ws.Range("B18").ClearContents()

# Update the selected cell (view state) to match the authored state
$ws.Range("I22").Select() | Out-Null
